# Update the Spotify export worksheet: rows 5-51 (Artist/Track/Album/Duration)
# get new values as described in the diff (one old row removed near the top,
# two old rows merged/removed in the middle, and three new rows appended at
# the bottom).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('What So Not', 'Lights Go Out', 'Lights Go Out', 4.01),
    @('BL3SS', 'R 2 ME', 'R 2 ME', 2.54),
    @('Barry Can''t Swim', 'Kimpton', 'Kimpton', 3.8),
    @('Daniel Allan', 'Something More', 'Something More', 2.7),
    @('Conrad.', 'we stayed up all night', 'we stayed up all night', 3.55),
    @('EVAN GIIA', 'MUD MELODIES', 'STAMINA', 3.17),
    @('ATTLAS', 'The One', 'The One', 3.26),
    @('gardenstate', 'The Best Part - Lost Prince Remix', 'The Best Part (Remixes)', 3.31),
    @('Eli Brown', 'Pulling Me Back', 'Pulling Me Back', 3.31),
    @('Tudor', 'Bonfire', 'Bonfire', 3.97),
    @('XANDRA', 'Feel Good - Punctual Remix', 'Feel Good (Punctual Remix)', 3.23),
    @('Return Of The Jaded', 'Be The Reason - Club Mix', 'Be The Reason (Club Mix)', 3.32),
    @('Punctual', 'You''re Not Alone - James Carter Remix', 'You''re Not Alone (James Carter Remix)', 3.03),
    @('HÜMAN', 'Craving You', 'Craving You', 4.45),
    @('Volyri', 'Better Now', 'Better Now', 2.59),
    @('HKLMR', 'Serenity (Kin Le Max Remix)', 'Serenity (Kin Le Max Remix)', 2.72),
    @('Calvin Harris', 'Blessings', 'Blessings', 3.66),
    @('KC Lights', 'CHOOSE LOVE', 'CHOOSE LOVE', 2.78),
    @('Spacey Jane', 'Weightless', 'Sunlight', 4.18),
    @('Spacey Jane', 'August', 'If That Makes Sense', 3.98),
    @('Spacey Jane', 'Ily the Most', 'If That Makes Sense', 2.83),
    @('Spacey Jane', 'Estimated Delivery', 'If That Makes Sense', 3.42),
    @('Spacey Jane', 'The More That it Hurts', 'If That Makes Sense', 3.03),
    @('Spacey Jane', 'So Much Taller', 'If That Makes Sense', 3.34),
    @('Spacey Jane', 'I Can’t Afford to Lose You', 'If That Makes Sense', 3.75),
    @('Spacey Jane', 'How to Kill Houseplants', 'If That Makes Sense', 3.48),
    @('Spacey Jane', 'Impossible to Say', 'If That Makes Sense', 3.53),
    @('Spacey Jane', 'All the Noise', 'If That Makes Sense', 3),
    @('Spacey Jane', 'Whateverrrr', 'If That Makes Sense', 2.97),
    @('Spacey Jane', 'Through My Teeth', 'If That Makes Sense', 3.42),
    @('Tchami', 'Praise', 'Year Zero', 3.5),
    @('Goodboys', 'Blindspot', 'Blindspot', 3.5),
    @('RAYE', 'Call On Me - KREAM Remix', 'Call On Me (KREAM Remix)', 3.68),
    @('Three Drives On A Vinyl', 'Greece 2000 - KREAM Remix', 'Greece 2000 (KREAM Remix)', 3.58),
    @('KREAM', 'Reverie', 'Reverie', 3.48),
    @('KREAM', 'Manta', 'Manta', 4.12),
    @('Koastle', 'Life (Can''t Get Much Better)', 'Life (Can''t Get Much Better)', 3.38),
    @('Effemar', 'Needing Space', 'Needing Space', 3.54),
    @('ATRIP', 'HERZSCHLAG', 'KLUBPARTEI', 3.9),
    @('HotLap', 'Set You Free', 'Set You Free', 3.4),
    @('Simon Doty', 'Tattoo', 'Tattoo', 3.65),
    @('Tina Says', 'Barriers', 'Barriers', 3.24),
    @('Tom Westy', 'Remember Me', 'Remember Me', 3.14),
    @('JEWELS', 'JETLAGGED', 'JETLAGGED', 4.55),
    @('Bad Friends', 'Will U?', 'Will U?', 2.71),
    @('Koastle', 'Dr. Phil', 'Dr. Phil', 3.39),
    @('Koastle', 'Sabotage', 'Sabotage', 3.27)
)

$startRow = 5
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
}
